# Updated spreadsheet & comments
# Adds a handful of "Tested" / "Mode" markers to the CANACC5 opcode table,
# then leaves the view scrolled/selected where the author left off.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "Tested" ticks (column E, centered style) -------------------------
$testedCells = "E27", "E31", "E32"
foreach ($addr in $testedCells) {
    $ws.Range($addr).Value = "y"
    $ws.Range($addr).HorizontalAlignment = -4108   # xlCenter
}

# --- New "Mode" notes (column D, default style) -----------------------------
$ws.Range("D30").Value = "setup mode"
$ws.Range("D34").Value = "learn mode"
$ws.Range("D37").Value = "setup mode"

# --- Restore the view position / selection the author ended up with --------
$ws.Range("D37").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 24 | Out-Null
$excel.ActiveWindow.ScrollColumn = 1 | Out-Null
